# apply hungarian algorithm for visualization
#
# Populates Sheet2 with the factorization/reconstruction-error comparison
# table (Kmeans / PCA / NMF) and makes Sheet2 the active/selected sheet,
# mirroring the authoring edit captured in the target diff.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- header row -----------------------------------------------------
$ws2.Range("A1").Value = "factorization"
$ws2.Range("B1").Value = "k"
$ws2.Range("C1").Value = "reconstruction error"
$ws2.Range("D1").Value = "norm(hyperfeat)"
$ws2.Range("E1").Value = "error/norm"

# --- Kmeans row -------------------------------------------------------
$ws2.Range("A2").Value = "Kmeans"
$ws2.Range("B2").Value = 4
$ws2.Range("C2").Value = 218.13
$ws2.Range("D2").Value = 333.15
$ws2.Range("E2").Formula = "=C2/D2"

# --- PCA row ------------------------------------------------------------
$ws2.Range("A3").Value = "PCA"
$ws2.Range("B3").Value = 4
$ws2.Range("C3").Value = 287.40118
$ws2.Range("D3").Value = 333.15
$ws2.Range("E3").Formula = "=C3/D3"

# --- NMF row --------------------------------------------------------
$ws2.Range("A4").Value = "NMF"
$ws2.Range("B4").Value = 4
$ws2.Range("C4").Value = 151.40488
$ws2.Range("D4").Value = 261.87
$ws2.Range("E4").Formula = "=C4/D4"

# --- column widths (characters) --------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 17.285714285714285
$ws2.Columns.Item(3).ColumnWidth = 20.142857142857142
$ws2.Columns.Item(4).ColumnWidth = 15.428571428571429
$ws2.Columns.Item(5).ColumnWidth = 11.857142857142858

# --- make Sheet2 the active sheet/tab, with K4 selected ---------------
$ws2.Activate() | Out-Null
$ws2.Range("K4").Select() | Out-Null
